$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 323, shifting the old row 323 (and everything
# below it, through the former row 437) down by one. This also extends the
# sheet's used range from A1:R437 to A1:R438.
$ws.Rows.Item(323).Insert()

# Copy the "static" columns from the row that used to be 323 (now 324) into
# the freshly inserted row 323, then set the new/changed values for this
# week's record (D, J, K, L, M, P).
$ws.Range("A323").Value = $ws.Range("A324").Value2
$ws.Range("B323").Value = $ws.Range("B324").Value2
$ws.Range("C323").Value = $ws.Range("C324").Value2
$ws.Range("E323").Value = $ws.Range("E324").Value2
$ws.Range("F323").Value = $ws.Range("F324").Value2
$ws.Range("G323").Value = $ws.Range("G324").Value2
$ws.Range("H323").Value = $ws.Range("H324").Value2
$ws.Range("I323").Value = $ws.Range("I324").Value2
$ws.Range("N323").Value = $ws.Range("N324").Value2
$ws.Range("O323").Value = $ws.Range("O324").Value2
$ws.Range("Q323").Value = $ws.Range("Q324").Value2
$ws.Range("R323").Value = $ws.Range("R324").Value2

$ws.Range("D323").Value = 44900
$ws.Range("J323").Value = 500
$ws.Range("K323").Value = 7000
$ws.Range("L323").Value = 8000
$ws.Range("M323").Value = 7500
$ws.Range("P323").Value = 375
